$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header labels in row 1 (columns D, E, F)
$ws.Range("D1").Value = "facility manager"
$ws.Range("E1").Value = "mobile number"
$ws.Range("F1").Value = "email id"

# Update the active selection on the sheet from E18 to J1
$ws.Range("J1").Select()
